# Update chapter numbering from "3.x"/"2.x" to "7.x" on the title shapes
# of slides 3, 4, 5 and 6 (titles: "3.1 Code" -> "7.1 Code",
# "2.1 Code" -> "7.1 Code", "3.2 Verify" -> "7.2 Verify" (x2)).

$p = $ppt.ActivePresentation

$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "7.1 Code"

$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "7.1 Code"

$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "7.2 Verify"

$s6 = $p.Slides.Item(6)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = "7.2 Verify"
